$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.282.33"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").Value = "3.505.64"
$ws.Range("E3").Value = "  -2.63%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'574.98"
$ws.Range("E5").Value = "  -1.01%  "

$ws.Range("D6").Value = "'185.68"
$ws.Range("E6").Value = "  -3.16%  "

$ws.Range("D7").Value = "3.496.37"
$ws.Range("E7").Value = "  -2.83%  "

$ws.Range("D8").Value = "'0.612"
$ws.Range("E8").Value = "  -3.34%  "

$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "'0.189"
$ws.Range("E10").Value = "  +4.02%  "

$ws.Range("E11").Value = "  -2.84%  "

$ws.Range("D12").Value = "'54.18"
$ws.Range("E12").Value = "  -3.17%  "

$ws.Range("D13").Value = "'0.0000301"
$ws.Range("E13").Value = "  -2.05%  "

$ws.Range("D14").Value = "'9.44"
$ws.Range("E14").Value = "  -2.82%  "

$ws.Range("D15").Value = "4.065.86"
$ws.Range("E15").Value = "  -2.83%  "

$ws.Range("D16").Value = "'19.34"
$ws.Range("E16").Value = "  -3.13%  "

$ws.Range("D17").Value = "69.252.12"
$ws.Range("E17").Value = "  -1.34%  "

$ws.Range("D18").Value = "3.502.54"
$ws.Range("E18").Value = "  -2.92%  "

$ws.Range("E19").Value = "  -3.25%  "

$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").Value = "'543.07"
$ws.Range("E21").Value = "  +12.59%  "

$ws.Range("E22").Value = "  -3.80%  "

$ws.Range("D23").Value = "'18.49"
$ws.Range("E23").Value = "  -4.17%  "

$ws.Range("E24").Value = "  -1.23%  "

$ws.Range("D25").Value = "'4.43"
$ws.Range("E25").Value = "  +0.43%  "

$ws.Range("D26").Value = "'93.93"
$ws.Range("E26").Value = "  -1.75%  "

$ws.Range("E27").Value = "  +1.47%  "

$ws.Range("D28").Value = "'2.94"
$ws.Range("E28").Value = "  -2.11%  "

$ws.Range("D29").Value = "'9.11"
$ws.Range("E29").Value = "  -3.00%  "

$ws.Range("D30").Value = "'31.83"
$ws.Range("E30").Value = "  -1.28%  "

$ws.Range("D31").Value = "'7.25"
$ws.Range("E31").Value = "  -6.37%  "

$ws.Range("D32").Value = "'12.59"
$ws.Range("E32").Value = "  +2.72%  "

$ws.Range("D33").Value = "'64.50"
$ws.Range("E33").Value = "  -3.35%  "

$ws.Range("E34").Value = "  -6.05%  "

$ws.Range("D35").Value = "'536.70"
$ws.Range("E35").Value = "  -8.75%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "'3.08"
$ws.Range("E36").Value = "  +8.05%  "

$ws.Range("B37").Value = "InjectiveProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D37").Value = "'37.94"
$ws.Range("E37").Value = "  -2.88%  "

$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.402"
$ws.Range("E38").Value = "  +0.92%  "

$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").Value = "0.0₃0766"
$ws.Range("E40").Value = "  -4.88%  "

$ws.Range("D41").Value = "'3.37"
$ws.Range("E41").Value = "  -2.85%  "

$ws.Range("E42").Value = "  -2.83%  "

$ws.Range("D43").Value = "3.310.17"
$ws.Range("E43").Value = "  +2.31%  "

$ws.Range("D44").Value = "'3.06"
$ws.Range("E44").Value = "  -8.17%  "

$ws.Range("E45").Value = "  -3.46%  "

$ws.Range("D46").Value = "'0.0444"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("E47").Value = "  +4.11%  "

$ws.Range("E49").Value = "  -6.63%  "

$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("D51").Value = "'137.43"
$ws.Range("E51").Value = "  +2.47%  "
